$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("B1").Value = "Buying Opportunity"
$ws.Range("C1").Value = "support Zone"
$ws.Range("D1").Value = "long buildup"
$ws.Range("E1").Value = "Short buildup"
$ws.Range("F1").Value = "FII ENTERING"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "NSE:AHLUCONT"
$ws.Range("C2").Value = "NSE:ACL"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:BATAINDIA"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "NSE:NIRAJ"
$ws.Range("C3").Value = "NSE:ANDHRAPAP"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:APOLLOTYRE"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = "NSE:ATGL"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:BALMLAWRIE"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:BANCOINDIA"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:BEPL"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:BPL"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:CARBORUNIV"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:CONSUMBEES"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:CYBERTECH"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:DELTAMAGNT"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DIAMONDYD"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:DIVGIITTS"
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:DPWIRES"
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:EKC"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:EMUDHRA"
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:EPIGRAL"
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:FACT"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "NSE:GABRIEL"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""

# Row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "NSE:GANESHHOUC"
$ws.Range("D22").Value = ""
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""

# Row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = "NSE:GENCON"
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = ""

# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "NSE:GMRP&UI"
$ws.Range("D24").Value = ""
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "NSE:GOLDIAM"
$ws.Range("D25").Value = ""
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""

# Row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "NSE:GSFC"
$ws.Range("D26").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""

# Row 27
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "NSE:GUJALKALI"
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""

# Row 28
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = "NSE:GULFOILLUB"
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""

# Row 29
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "NSE:GULPOLY"
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""

# Row 30
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = "NSE:HINDMOTORS"
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = ""
$ws.Range("F30").Value = ""

# Row 31
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = "NSE:INDOAMIN"
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = ""

# Row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = "NSE:INFIBEAM"
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = ""

# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = "NSE:ITDC"
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = ""

# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = ""
$ws.Range("C34").Value = "NSE:JAYSREETEA"
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = ""

# Row 35
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = ""
$ws.Range("C35").Value = "NSE:KOTARISUG"
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = ""
$ws.Range("F35").Value = ""

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = ""
$ws.Range("C36").Value = "NSE:LUMAXIND"
$ws.Range("D36").Value = ""
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = ""

# Row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = ""
$ws.Range("C37").Value = "NSE:LXCHEM"
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = ""
$ws.Range("F37").Value = ""

# Row 38
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = ""
$ws.Range("C38").Value = "NSE:MAGNUM"
$ws.Range("D38").Value = ""
$ws.Range("E38").Value = ""
$ws.Range("F38").Value = ""

# Row 39
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = ""
$ws.Range("C39").Value = "NSE:MAXESTATES"
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = ""
$ws.Range("F39").Value = ""

# Row 40
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = ""
$ws.Range("C40").Value = "NSE:MGEL"
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = ""
$ws.Range("F40").Value = ""

# Row 41
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = ""
$ws.Range("C41").Value = "NSE:MICEL"
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = ""
$ws.Range("F41").Value = ""

# Row 42
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = ""
$ws.Range("C42").Value = "NSE:MOLDTECH"
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = ""
$ws.Range("F42").Value = ""

# Row 43
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = ""
$ws.Range("C43").Value = "NSE:MONARCH"
$ws.Range("D43").Value = ""
$ws.Range("E43").Value = ""
$ws.Range("F43").Value = ""

# Row 44
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = ""
$ws.Range("C44").Value = "NSE:MRF"
$ws.Range("D44").Value = ""
$ws.Range("E44").Value = ""
$ws.Range("F44").Value = ""

# Row 45
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = ""
$ws.Range("C45").Value = "NSE:NAVINIFTY"
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = ""
$ws.Range("F45").Value = ""

# Row 46
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = ""
$ws.Range("C46").Value = "NSE:NAVNETEDUL"
$ws.Range("D46").Value = ""
$ws.Range("E46").Value = ""
$ws.Range("F46").Value = ""

# Row 47
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = ""
$ws.Range("C47").Value = "NSE:NIITLTD"
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = ""
$ws.Range("F47").Value = ""

# Row 48
$ws.Range("A2").Copy($ws.Range("A48"))
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = "NSE:NLCINDIA"
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = ""
$ws.Range("F48").Value = ""

# Row 49
$ws.Range("A2").Copy($ws.Range("A49"))
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = ""
$ws.Range("C49").Value = "NSE:OCCL"
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("F49").Value = ""

# Row 50
$ws.Range("A2").Copy($ws.Range("A50"))
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = ""
$ws.Range("C50").Value = "NSE:PATINTLOG"
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = ""
$ws.Range("F50").Value = ""

# Row 51
$ws.Range("A2").Copy($ws.Range("A51"))
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = ""
$ws.Range("C51").Value = "NSE:PRICOLLTD"
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = ""
$ws.Range("F51").Value = ""

# Row 52
$ws.Range("A2").Copy($ws.Range("A52"))
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = ""
$ws.Range("C52").Value = "NSE:PSPPROJECT"
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = ""

# Row 53
$ws.Range("A2").Copy($ws.Range("A53"))
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = ""
$ws.Range("C53").Value = "NSE:RADIANTCMS"
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = ""

# Row 54
$ws.Range("A2").Copy($ws.Range("A54"))
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = ""
$ws.Range("C54").Value = "NSE:REPRO"
$ws.Range("D54").Value = ""
$ws.Range("E54").Value = ""
$ws.Range("F54").Value = ""

# Row 55
$ws.Range("A2").Copy($ws.Range("A55"))
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = ""
$ws.Range("C55").Value = "NSE:RIIL"
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = ""

# Row 56
$ws.Range("A2").Copy($ws.Range("A56"))
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = ""
$ws.Range("C56").Value = "NSE:RKEC"
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("F56").Value = ""

